# Full refreshed data set for the 'ランサーズ' sheet (rows 2-16).
# Columns: A=取得日時 B=タイトル C=カテゴリ D=価格 E=締切 F=URL G=優先度スコア H=スキル概要
$data = @(
  @("2025-10-01 01:24:41", "【AI活用】データ分析Webサービス開発パートナー募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5404026", "368", "🔥AI,Ai ◆開発"),
  @("2025-10-01 01:24:41", "【限定タスク】SIM AI の Google 認証ログイン機能の「最終調整」のみ代行(環境構築済み)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403583", "298", "🔥AI,Ai"),
  @("2025-10-01 01:24:41", "【急募】メモリデータ管理ツール開発のプロフェッショナル募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5393508", "158", "◆ツール,開発 ◇管理"),
  @("2025-10-01 01:24:41", "システム開発において活躍できる案件紹介", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397117", "125", "◆開発,システム開発"),
  @("2025-10-01 01:24:41", "検索エンジン予測変換(サジェスト表示)に関するツール開発・調査依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403988", "123", "◆ツール,開発"),
  @("2025-10-01 01:24:41", "【3万円/Webツール開発】サジェスト対策", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403789", "123", "◆ツール,開発"),
  @("2025-10-01 01:24:41", "【RPA/Power Automate】税務システム自動化プロジェクトの依頼", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403634", "103", "◆自動化"),
  @("2025-10-01 01:24:41", "急募バックエンドエンジニア マッチングサイトの開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5404059", "93", "◆開発 ◇サイト"),
  @("2025-10-01 01:24:41", "【急募】LINE WORKSで定期メッセージ配信ツール作成依頼", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403166", "65", "◆ツール"),
  @("2025-10-01 01:24:41", "【急募】教育系のWEBサイトの作成", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403527", "33", "◇サイト"),
  @("2025-10-01 01:24:41", "要件定義や基本設計ができる方(1人月、約2年アサイン予定)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5391221", "25", ""),
  @("2025-10-01 01:24:41", "【ビジネスパートナー募集】英語が得意な方を探しています", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403384", "25", ""),
  @("2025-10-01 01:24:41", "初回 サブスクペイからCSVデータをダウンロードし、データベース同期するプログラムの作成", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5403072", "18", ""),
  @("2025-10-01 01:24:41", "限定公開 PR 限定公開の仕事", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399347", "13", ""),
  @("2025-10-01 01:24:41", "急募 【占い師向け】VBAマクロで自動反映エクセル作成依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5404010", "13", ""),
)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Drop every existing hyperlink (relationships + <hyperlink> entries) so the
# F-column links can be rebuilt from scratch in the new row order; this
# leaves each cell's "Hyperlink" style (s="1") untouched.
$ws.Cells.Hyperlinks.Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2

  $ws.Cells.Item($r, 1).Value2 = $row[0]          # A: 取得日時
  $ws.Cells.Item($r, 2).Value2 = $row[1]          # B: タイトル
  $ws.Cells.Item($r, 3).Value2 = $row[2]          # C: カテゴリ
  $ws.Cells.Item($r, 4).Value2 = $row[3]          # D: 価格
  $ws.Cells.Item($r, 5).Value2 = $row[4]          # E: 締切

  $fCell = $ws.Cells.Item($r, 6)                  # F: URL (hyperlink)
  $fCell.Value2 = $row[5]
  $ws.Hyperlinks.Add($fCell, $row[5])
  $fCell.Style = "Hyperlink"

  $ws.Cells.Item($r, 7).Value2 = [int]$row[6]     # G: 優先度スコア (numeric)

  if ($row[7] -ne "") {
    $ws.Cells.Item($r, 8).Value2 = $row[7]        # H: スキル概要 (omitted when blank)
  }
}
